{"js": "// The author repositioned the cursor (and thus Word's \"_GoBack\" bookmark,\n// which always tracks the location of the last edit) from the end of the\n// \"Data analysis\" paragraph to a point in the middle of the BCI-imagery\n// sentence, right after \"...motor imagery signals can be detected\" and\n// before \" using EEG signals to help people...\". No visible text changed;\n// the run was simply split at that point and the bookmark moved there.\n\n// 1. Remove the existing \"_GoBack\" bookmark (a document can only have one\n//    bookmark with a given name).\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2. Find the unique anchor text and split the run right after \"detected\".\nconst body = context.document.body;\nconst results = body.search(\"can be detected\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const match = results.items[0];\n  const insertionPoint = match.getRange(\"End\");\n\n  // 3. Re-insert \"_GoBack\" collapsed at that exact spot.\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The author repositioned the cursor (and thus Word's \"_GoBack\" bookmark,\n# which always tracks the location of the last edit) from the end of the\n# \"Data analysis\" paragraph to a point in the middle of the BCI-imagery\n# sentence, right after \"...motor imagery signals can be detected\" and\n# before \" using EEG signals to help people...\". No visible text changed;\n# the run was simply split at that point and the bookmark moved there.\n\n$d = $word.ActiveDocument\n\n# 1. Remove the existing \"_GoBack\" bookmark (a document can only have one\n#    bookmark with a given name).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2. Find the unique anchor text and collapse to the point right after it.\n$r = $d.Content\n$r.Find.Text = \"can be detected\"\n$r.Find.Execute() | Out-Null\n\nif ($r.Find.Found) {\n    $r.Collapse(0)  # wdCollapseEnd\n\n    # 3. Re-insert \"_GoBack\" collapsed at that exact spot.\n    $d.Bookmarks.Add(\"_GoBack\", $r)\n}\n"}
